$d = $word.ActiveDocument
$d.Content.Find.Execute("826×3=2478", $true, $false, $false, $false, $false, $true, 1, $false, "734×7=5138", 2)
$d.Content.Find.Execute("259×8=2072", $true, $false, $false, $false, $false, $true, 1, $false, "743×6=4458", 2)
$d.Content.Find.Execute("206×8=1648", $true, $false, $false, $false, $false, $true, 1, $false, "343×5=1715", 2)
$d.Content.Find.Execute("982×6=5892", $true, $false, $false, $false, $false, $true, 1, $false, "623×5=3115", 2)
$d.Content.Find.Execute("319×3=957", $true, $false, $false, $false, $false, $true, 1, $false, "558×7=3906", 2)
$d.Content.Find.Execute("837×4=3348", $true, $false, $false, $false, $false, $true, 1, $false, "401×4=1604", 2)
$d.Content.Find.Execute("433×6=2598", $true, $false, $false, $false, $false, $true, 1, $false, "528×7=3696", 2)
$d.Content.Find.Execute("909×3=2727", $true, $false, $false, $false, $false, $true, 1, $false, "457×8=3656", 2)
$d.Content.Find.Execute("626×4=2504", $true, $false, $false, $false, $false, $true, 1, $false, "329×7=2303", 2)
$d.Content.Find.Execute("677×3=2031", $true, $false, $false, $false, $false, $true, 1, $false, "234×6=1404", 2)
$d.Content.Find.Execute("278×3=834", $true, $false, $false, $false, $false, $true, 1, $false, "301×5=1505", 2)
$d.Content.Find.Execute("236×2=472", $true, $false, $false, $false, $false, $true, 1, $false, "826×8=6608", 2)
$d.Content.Find.Execute("676×8=5408", $true, $false, $false, $false, $false, $true, 1, $false, "386×9=3474", 2)
$d.Content.Find.Execute("542×4=2168", $true, $false, $false, $false, $false, $true, 1, $false, "702×7=4914", 2)
$d.Content.Find.Execute("543×2=1086", $true, $false, $false, $false, $false, $true, 1, $false, "448×4=1792", 2)
$d.Content.Find.Execute("473×4=1892", $true, $false, $false, $false, $false, $true, 1, $false, "713×6=4278", 2)
$d.Content.Find.Execute("778×9=7002", $true, $false, $false, $false, $false, $true, 1, $false, "868×6=5208", 2)
$d.Content.Find.Execute("244×9=2196", $true, $false, $false, $false, $false, $true, 1, $false, "710×6=4260", 2)
$d.Content.Find.Execute("344×2=688", $true, $false, $false, $false, $false, $true, 1, $false, "458×4=1832", 2)
$d.Content.Find.Execute("203×4=812", $true, $false, $false, $false, $false, $true, 1, $false, "911×9=8199", 2)
$d.Content.Find.Execute("552×2=1104", $true, $false, $false, $false, $false, $true, 1, $false, "432×9=3888", 2)
$d.Content.Find.Execute("724×5=3620", $true, $false, $false, $false, $false, $true, 1, $false, "227×8=1816", 2)
$d.Content.Find.Execute("575×4=2300", $true, $false, $false, $false, $false, $true, 1, $false, "746×2=1492", 2)
$d.Content.Find.Execute("829×2=1658", $true, $false, $false, $false, $false, $true, 1, $false, "257×3=771", 2)
$d.Content.Find.Execute("507×9=4563", $true, $false, $false, $false, $false, $true, 1, $false, "980×2=1960", 2)
